$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows hold per-innings stats (runs, balls, fours) for the same player/team
# and the diff shows the five data rows (2-6) being cyclically rotated:
#   new Row2 = old Row6 (26, 21, 3)
#   new Row3 = old Row2 (7, 12, 1)
#   new Row4 = old Row5 (3, 10, 0)
#   new Row5 = old Row3 (22, 16, 3)
#   new Row6 = old Row4 (4, 7, 0)
# Columns A, B, F are unchanged across all rows, so only C (runs), D (balls)
# and E (fours) need to be rewritten with their new values.

$ws.Range("C2").Value = "26"
$ws.Range("D2").Value = "21"
$ws.Range("E2").Value = "3"

$ws.Range("C3").Value = "7"
$ws.Range("D3").Value = "12"
$ws.Range("E3").Value = "1"

$ws.Range("C4").Value = "3"
$ws.Range("D4").Value = "10"
$ws.Range("E4").Value = "0"

$ws.Range("C5").Value = "22"
$ws.Range("D5").Value = "16"
$ws.Range("E5").Value = "3"

$ws.Range("C6").Value = "4"
$ws.Range("D6").Value = "7"
$ws.Range("E6").Value = "0"
